$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.268.83"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.791.43"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'432.69"
$ws.Range("E5").Value = "  +5.47%  "
$ws.Range("D6").Value = "'138.39"
$ws.Range("E6").Value = "  +4.05%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'0.152"
$ws.Range("E10").Value = "  -8.56%  "
$ws.Range("D11").Value = "'0.0000313"
$ws.Range("E11").Value = "  -13.14%  "
$ws.Range("D12").Value = "'42.67"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").Value = "'10.44"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "4.393.41"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "'15.00"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.773.06"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'19.98"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("D20").Value = "66.375.35"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'406.17"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'14.87"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").Value = "'3.26"
$ws.Range("E23").Value = "  +6.10%  "
$ws.Range("D24").Value = "'84.87"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'36.77"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'3.33"
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.81"
$ws.Range("E27").Value = "  +32.60%  "
$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +4.98%  "
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").Value = "'0.139"
$ws.Range("E30").Value = "  +14.20%  "
$ws.Range("D31").Value = "'13.80"
$ws.Range("E31").Value = "  +11.44%  "
$ws.Range("D32").Value = "'706.68"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'41.73"
$ws.Range("E34").Value = "  +5.87%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'5.66"
$ws.Range("E36").Value = "  +32.52%  "
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "'56.19"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +39.25%  "
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").Value = "'0.142"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("D43").Value = "0.0₃0679"
$ws.Range("E43").Value = "  -7.70%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "'0.332"
$ws.Range("E45").Value = "  +11.68%  "
$ws.Range("D46").Value = "'3.26"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").Value = "'2.07"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "'140.24"
$ws.Range("E50").Value = "  -3.96%  "
$ws.Range("E51").Value = "  -0.11%  "

# Reset style on text-forced numeric-looking price cells to avoid stray quotePrefix formatting
$textForcedCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D15", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D36", "D38", "D40", "D41", "D42", "D44", "D45", "D46", "D49", "D50")
foreach ($ref in $textForcedCells) {
    $ws.Range($ref).Style = "Normal"
}
